# Replace the "word" column (B2:B193) with the new working set of sequences.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "gleichen", "fischen", "heilen", "schaden", "tropfen", "flüchten", "drücken", "retten",
    "sorgen", "saufen", "schalten", "formen", "stoppen", "heulen", "schütteln", "brauchen",
    "bauen", "arten", "machen", "weinen", "dauern", "fahren", "rufen", "pflegen",
    "bergen", "schreiten", "kosten", "hupen", "rücken", "wachsen", "eignen", "scheinen",
    "rauben", "erben", "spüren", "fassen", "trauen", "treten", "wenden", "plaudern",
    "schlafen", "segeln", "mauern", "sperren", "räumen", "herzen", "steuern", "hauen",
    "beten", "zeigen", "wundern", "wüten", "baden", "scheitern", "atmen", "siegen",
    "lügen", "pfeifen", "leisten", "sinken", "gelten", "schmecken", "weigern", "spielen",
    "trennen", "schicken", "wehtun", "dienen", "zielen", "schwächen", "wehen", "reizen",
    "ärgern", "grüßen", "wiegen", "backen", "geben", "trösten", "lesen", "graben",
    "zünden", "küssen", "flehen", "narren", "dringen", "schenken", "stecken", "klingen",
    "füllen", "brauen", "sprengen", "boxen", "malen", "loben", "handeln", "betteln",
    "helfen", "schwingen", "kehren", "ändern", "danken", "tragen", "fallen", "fließen",
    "mühen", "spinnen", "schwören", "äußern", "greifen", "ächzen", "pflanzen", "sichern",
    "münzen", "decken", "biegen", "drucken", "liefern", "testen", "gründen", "folgen",
    "weichen", "ehren", "jagen", "tauschen", "treiben", "wetten", "träumen", "schrecken",
    "jubeln", "werfen", "runden", "bitten", "bellen", "starren", "liegen", "bleiben",
    "fangen", "suchen", "tanzen", "rasen", "mögen", "zögern", "platzen", "stehlen",
    "teilen", "streichen", "freuen", "führen", "enden", "schulden", "filmen", "läuten",
    "wohnen", "wirken", "fällen", "hören", "sterben", "heben", "süßen", "segnen",
    "kaufen", "quälen", "kümmern", "drehen", "klettern", "stoßen", "töten", "warnen",
    "kichern", "schlucken", "knarren", "duschen", "wechseln", "klagen", "knien", "werden",
    "hindern", "sagen", "planen", "tollen", "irren", "lockern", "hassen", "altern",
    "rühren", "bluten", "feiern", "stammen", "prügeln", "achten", "seufzen", "lieben"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}
